# model and vol 26.02.18
# Add a new data row (row 40) to the MESAlteckWater sheet, mirroring the
# pattern of the previous rows (A: AFP-W{n}-1, B: AO-W{n}-1, C: W{n}-1,
# D/E/F: model + volume values), and move the active selection to G40.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MESAlteckWater")

# Row 39 (previous last row) picked up a volume reading of 0 too.
$ws.Range("E39").Value = 0

# New row values
$ws.Range("A40").Value = "AFP-W16-1"
$ws.Range("B40").Value = "AO-W16-1"
$ws.Range("C40").Value = "W16-1"
$ws.Range("D40").Value = 55
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0.0056000000000011596

# Carry over the same cell formatting used by the row above (A39/B39/C39),
# and the plain centered style used throughout column D, so the new row
# visually matches the rest of the table.
$ws.Range("A39").Copy() | Out-Null
$ws.Range("A40").PasteSpecial(-4122) | Out-Null

$ws.Range("B39:C39").Copy() | Out-Null
$ws.Range("B40:C40").PasteSpecial(-4122) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move/select the active cell like the author left it (G40).
$ws.Range("G40").Select() | Out-Null
